$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = @{
    2 = 2015
    3 = 2016
    4 = 2017
    5 = 2018
    6 = 2019
    7 = 2020
    8 = 2021
    9 = 2022
    10 = 2023
}

foreach ($row in $years.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $years[$row]
    $cell.NumberFormat = "General"
}

$ws.Columns.Item(1).ColumnWidth = 10.7109375
